$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 2: fix "wgde" + "v" (split across a stray _GoBack bookmark)
# into a single contiguous run "wgdev". Using Find/Replace merges the
# two runs into one and removes the bookmark that sat between them
# (matching real Word's "replace consumes enclosed bookmark" semantics).
# ------------------------------------------------------------------
$rngWgdev = $d.Content
$rngWgdev.Find.Execute("wgdev", $true, $false, $false, $false, $false, $true, 1, $false, "wgdev", 2)

# ------------------------------------------------------------------
# Change 1: split "beta_SE" (in "Add beta_SE to Dataset") into
# "beta_S" + "E", re-typing just the trailing "E" and dropping a new
# _GoBack bookmark exactly at the split point (this is where Word's
# last-edit bookmark ends up after the most recent keystroke).
# ------------------------------------------------------------------
$text = $d.Content.Text
$idx = $text.IndexOf("beta_SE")

$rngOldE = $d.Range($idx + 6, $idx + 7)
$rngOldE.Delete()
$rngNewE = $d.Range($idx + 6, $idx + 6)
$rngNewE.InsertAfter("E")

$rngMark = $d.Range($idx + 6, $idx + 6)
$d.Bookmarks.Add("_GoBack", $rngMark)

# ------------------------------------------------------------------
# Change 3: mark the "Update help, test, vignette" bullet as fully
# crossed-out (strike the whole run of text plus the paragraph mark),
# then normalize it back into a single run.
# ------------------------------------------------------------------
$text2 = $d.Content.Text
$idx2 = $text2.IndexOf("Update help, test, vignette")
$target = "Update help, test, vignette"

$rngPara = $d.Range($idx2, $idx2 + $target.Length)
$para = $rngPara.Paragraphs(1)
$para.Range.Font.StrikeThrough = 1

$rngMerge = $d.Content
$rngMerge.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, $target, 2)
